$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.667.27"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.447.00"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.89"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.53"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  +13.19%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.448.83"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").Value = "4.048.31"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000190"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.11"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "64.764.75"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "3.465.24"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.26"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.14"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.09"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").Value = "  +3.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.16"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.94"
$ws.Range("E27").Value = "  +6.31%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +9.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.47"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.09"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +11.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.07"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  +5.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0772"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "2.956.49"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.40"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.59"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.53"
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.74"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0317"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.45"
$ws.Range("E46").Value = "  +10.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.09"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("E48").Value = "  +8.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "310.49"
$ws.Range("E49").Value = "  +6.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.60"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.859"
$ws.Range("E51").Value = "  +3.65%  "
